$d = $word.ActiveDocument

# The "Requisitos" list is a single ListBullet paragraph made of 21 runs,
# each run holding "<code> -  <name>  (Requisito)" followed by a line
# break. The edit re-sorts those 21 lines alphabetically by course code,
# keeping the exact same wording/run-per-line structure.

$newTexts = @(
  "LOB1003 -  Cálculo I  (Requisito)",
  "LOB1004 -  Cálculo II  (Requisito)",
  "LOB1006 -  Cálculo IV  (Requisito)",
  "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)",
  "LOB1011 -  Eletricidade Aplicada  (Requisito)",
  "LOB1018 -  Física I  (Requisito)",
  "LOB1019 -  Física II  (Requisito)",
  "LOB1024 -  Mecânica  (Requisito)",
  "LOB1036 -  Geometria Analítica  (Requisito)",
  "LOB1037 -  Àlgebra Linear  (Requisito)",
  "LOB1038 -  Física Experimental I  (Requisito)",
  "LOB1039 -  Física Experimental III  (Requisito)",
  "LOB1040 -  Laboratório de Eletricidade  (Requisito)",
  "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)",
  "LOB1052 -  Cálculo III  (Requisito)",
  "LOB1053 -  Física III  (Requisito)",
  "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)",
  "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
  "LOQ4095 -  Química Geral Experimental  (Requisito)",
  "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)",
  "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)"
)

# Locate the "Requisitos" heading paragraph; the bullet list is the very
# next paragraph.
$reqHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Requisitos") {
        $reqHeadingIndex = $i
        break
    }
}

$listPara = $d.Paragraphs($reqHeadingIndex + 1)
$listRange = $listPara.Range

# Build the XML-escaped <w:r> chunk for every line (text run + line break),
# matching the original run-per-line layout.
function Escape-Xml([string]$s) {
    $s = $s -replace "&", "&amp;"
    $s = $s -replace "<", "&lt;"
    $s = $s -replace ">", "&gt;"
    return $s
}

$runsXml = ""
foreach ($t in $newTexts) {
    $runsXml += "<w:r><w:t>" + (Escape-Xml $t) + "</w:t><w:br/></w:r>"
}

$packageXml = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
  "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
  "<pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
  "<w:body><w:p>" + $runsXml + "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# Remove the existing 21 runs (keep the paragraph mark / ListBullet style)
# then inject the freshly ordered runs via InsertXML so each line stays
# its own <w:r> instead of being coalesced into one run.
$deleteRange = $d.Range($listRange.Start, $listRange.End - 1)
if ($deleteRange.Start -lt $deleteRange.End) {
    $deleteRange.Delete()
}

$listPara2 = $d.Paragraphs($reqHeadingIndex + 1)
$insPos = $listPara2.Range.Start
$insertPoint = $d.Range($insPos, $insPos)
$insertPoint.InsertXML($packageXml)

Write-Host "Requisitos list re-sorted."
